$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 641, shifting existing rows 641+ down by one
# (old row 641 "2026/12/29" becomes new row 642, etc.).
$ws.Rows.Item(641).Insert()

# Populate the newly inserted row with the new data point.
# Force the date column to text formatting first so Excel does not
# reinterpret "2026/01/13" as a date serial number.
$ws.Cells.Item(641, 1).NumberFormat = "@"
$ws.Cells.Item(641, 1).Value = "2026/01/13"
$ws.Cells.Item(641, 2).Value = "火"
$ws.Cells.Item(641, 3).Value = 13
$ws.Cells.Item(641, 4).Value = 201
